{"js": "// Move the \"LOT2053 - Microbiologia (Requisito fraco)\" requisite line so it\n// appears FIRST in the \"Requisitos\" bullet list, ahead of the \"LOT2007 -\n// Bioqu\u00edmica I\" and \"LOT2040 - Engenharia Gen\u00e9tica\" lines (which keep their\n// relative order).\nconst MOVE_LINE = \"LOT2053 -  Microbiologia  (Requisito fraco)\";\nconst LINE_BREAK = \"\\u000b\"; // Office.js text representation of <w:br/>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the bullet-list paragraph that holds the requirement lines (the one\n// whose text contains our target line alongside the others).\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(MOVE_LINE) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not locate the paragraph containing \"' + MOVE_LINE + '\"');\n}\n\n// Already first? Nothing to do.\nif (!target.text.startsWith(MOVE_LINE)) {\n  // 1) Prepend the moved line (with its trailing break) at the very start of\n  //    the paragraph \u2014 this mints a new run and leaves the existing runs\n  //    (and their own breaks) completely untouched.\n  target.getRange(\"Start\").insertText(MOVE_LINE + LINE_BREAK, \"Before\");\n  await context.sync();\n\n  // 2) Remove the original occurrence of that line (the one that is now\n  //    duplicated further down in the paragraph, still followed by its\n  //    break) by searching for the exact \"text + break\" span and deleting\n  //    the LAST match (the first match is the copy we just inserted).\n  const hits = target.search(MOVE_LINE + LINE_BREAK, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  hits.items[hits.items.length - 1].delete();\n  await context.sync();\n}\n", "ps1": "# Move the \"LOT2053 - Microbiologia (Requisito fraco)\" requirement line so it\n# becomes the FIRST line of the \"Requisitos\" bullet list, ahead of the\n# \"LOT2007 - Bioquimica I\" and \"LOT2040 - Engenharia Genetica\" lines (which\n# keep their existing relative order).\n\n$d = $word.ActiveDocument\n\n$moveLine = \"LOT2053 -  Microbiologia  (Requisito fraco)\"\n$lineBreak = [char]11   # Word's manual line break (<w:br/>) as plain text\n\n# Locate the bullet-list paragraph that contains all three requirement lines.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($moveLine)) {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not locate the paragraph containing '$moveLine'\"\n}\n\nif (-not $target.Range.Text.StartsWith($moveLine)) {\n    # 1) Prepend the moved line (with its trailing break) at the very start\n    #    of the paragraph. This mints a brand-new run and leaves the\n    #    existing runs (and their own breaks) completely untouched.\n    $insertionPoint = $target.Range.Duplicate\n    $insertionPoint.Collapse(1)   # wdCollapseStart\n    $insertionPoint.InsertBefore($moveLine + $lineBreak)\n\n    # 2) Remove the original occurrence of that line (now duplicated further\n    #    down in the paragraph, still followed by its own break). Walk all\n    #    matches of \"text + break\" inside the paragraph and delete the LAST\n    #    one (the first match is the copy we just inserted at the start).\n    $needle = $moveLine + $lineBreak\n    $pStart = $target.Range.Start\n    $pEnd = $target.Range.End\n\n    $matches = @()\n    $cursor = $d.Range($pStart, $pEnd)\n    while ($true) {\n        $found = $cursor.Find.Execute($needle)\n        if (-not $found) { break }\n        if ($cursor.Start -ge $pEnd) { break }\n        $matches += , @($cursor.Start, $cursor.End)\n        $nextStart = $cursor.End\n        if ($nextStart -ge $pEnd) { break }\n        $cursor = $d.Range($nextStart, $pEnd)\n    }\n\n    if ($matches.Count -gt 0) {\n        $last = $matches[$matches.Count - 1]\n        $dupRange = $d.Range($last[0], $last[1])\n        $dupRange.Delete()\n    }\n}\n"}
